# Fruta / hortaliza, semanal
#
# A new weekly observation was appended at the top of the "Plátano"
# series (now dated 2021-09-08 / serial 44447). Every existing record
# from the previous row 73 to the bottom (previously ending at row 129)
# shifts down by one row, and the record that used to live in row 72
# (dated 2021-03-12 / serial 44267) is preserved verbatim as the new
# row 73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the current row 72 (the most-recent record before the edit)
# onto the clipboard, then insert a new blank row at 73, pushing
# everything below down by one.
$ws.Range("A72:T72").Copy()
$ws.Rows.Item(73).Insert()

# Paste the duplicated record into the freshly inserted row 73 — this
# is what used to be row 72, now preserved one row lower.
$ws.Range("A73:T73").PasteSpecial()

# Row 72 itself now represents the brand-new weekly observation; only
# its date actually changes (2021-03-12 -> 2021-09-08), every other
# attribute stays the same as it already was.
$newDate = Get-Date -Year 2021 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("D72").Value = $newDate

Write-Output "done"
